# Updates cached market-price / profit figures on the Leve profit sheets.
# Source data refreshed by the scheduled runner; only numeric columns H:N change.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H9").Value = 234.44444
$ws.Range("I9").Value = 31
$ws.Range("J9").Value = 259.875
$ws.Range("K9").Value = 31
$ws.Range("L9").Value = 259.875
$ws.Range("M9").Value = 138
$ws.Range("N9").Value = -597.875

$ws.Range("H28").Value = 1691.8695
$ws.Range("I28").Value = 2549.1538
$ws.Range("J28").Value = 577.4
$ws.Range("K28").Value = 2549.1538
$ws.Range("L28").Value = 577.4
$ws.Range("M28").Value = -2064.1538
$ws.Range("N28").Value = -1547.4

$ws.Range("H98").Value = 7416733
$ws.Range("I98").Value = 10301.077
$ws.Range("J98").Value = 55558540
$ws.Range("K98").Value = 10301.077
$ws.Range("L98").Value = 55558540
$ws.Range("M98").Value = -8803.076999999999
$ws.Range("N98").Value = -55561536

$ws.Range("H122").Value = 7416733
$ws.Range("I122").Value = 10301.077
$ws.Range("J122").Value = 55558540
$ws.Range("K122").Value = 30903.231
$ws.Range("L122").Value = 166675620
$ws.Range("M122").Value = -28453.231
$ws.Range("N122").Value = -166680520

$ws.Range("H132").Value = 1152.8334
$ws.Range("I132").Value = 1092.0785
$ws.Range("J132").Value = 2185.6667
$ws.Range("K132").Value = 3276.2355
$ws.Range("L132").Value = 6557.000100000001
$ws.Range("M132").Value = -746.2355000000002
$ws.Range("N132").Value = -11617.0001

$ws.Range("H138").Value = 4207.2964
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4207.2964
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 12621.8892
$ws.Range("N138").Value = -22901.8892
$ws.Range("M138").ClearContents()


$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 4649.58
$ws.Range("I32").Value = 4244.6045
$ws.Range("J32").Value = 8744.333000000001
$ws.Range("K32").Value = 4244.6045
$ws.Range("L32").Value = 8744.333000000001
$ws.Range("M32").Value = -3957.6045
$ws.Range("N32").Value = -9318.333000000001

$ws.Range("H61").Value = 3881.0527
$ws.Range("I61").Value = 3928.4666
$ws.Range("J61").Value = 3703.25
$ws.Range("K61").Value = 3928.4666
$ws.Range("L61").Value = 3703.25
$ws.Range("M61").Value = -3716.4666
$ws.Range("N61").Value = -4127.25

$ws.Range("H74").Value = 923.1429000000001
$ws.Range("I74").Value = 658.2222
$ws.Range("J74").Value = 1400
$ws.Range("K74").Value = 658.2222
$ws.Range("L74").Value = 1400
$ws.Range("M74").Value = 215.7778
$ws.Range("N74").Value = -3148

$ws.Range("H77").Value = 923.1429000000001
$ws.Range("I77").Value = 658.2222
$ws.Range("J77").Value = 1400
$ws.Range("K77").Value = 3291.111
$ws.Range("L77").Value = 7000
$ws.Range("M77").Value = 1076.889
$ws.Range("N77").Value = -15736

$ws.Range("H136").Value = 3881.0527
$ws.Range("I136").Value = 3928.4666
$ws.Range("J136").Value = 3703.25
$ws.Range("K136").Value = 11785.3998
$ws.Range("L136").Value = 11109.75
$ws.Range("M136").Value = -9235.399800000001
$ws.Range("N136").Value = -16209.75


$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H105").Value = 1722.6086
$ws.Range("I105").Value = 1730.5883
$ws.Range("J105").Value = 1700
$ws.Range("K105").Value = 1730.5883
$ws.Range("L105").Value = 1700
$ws.Range("M105").Value = 16.41170000000011
$ws.Range("N105").Value = -5194


$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H62").Value = 3990.238
$ws.Range("I62").Value = 3999.6667
$ws.Range("J62").Value = 3966.6667
$ws.Range("K62").Value = 3999.6667
$ws.Range("L62").Value = 3966.6667
$ws.Range("M62").Value = -3375.6667
$ws.Range("N62").Value = -5214.6667

$ws.Range("H65").Value = 3990.238
$ws.Range("I65").Value = 3999.6667
$ws.Range("J65").Value = 3966.6667
$ws.Range("K65").Value = 19998.3335
$ws.Range("L65").Value = 19833.3335
$ws.Range("M65").Value = -16878.3335
$ws.Range("N65").Value = -26073.3335


$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H2").Value = 212.46428
$ws.Range("I2").Value = 73.38461
$ws.Range("J2").Value = 333
$ws.Range("K2").Value = 440.3076599999999
$ws.Range("L2").Value = 1998
$ws.Range("M2").Value = -327.3076599999999
$ws.Range("N2").Value = -2224

$ws.Range("H70").Value = 1230
$ws.Range("I70").Value = 768
$ws.Range("K70").Value = 2304
$ws.Range("M70").Value = -1989

$ws.Range("H73").Value = 1230
$ws.Range("I73").Value = 768
$ws.Range("K73").Value = 2304
$ws.Range("M73").Value = -1212

$ws.Range("H102").Value = 4116.6665
$ws.Range("J102").Value = 4780
$ws.Range("L102").Value = 14340
$ws.Range("N102").Value = -19208

$ws.Range("H103").Value = 733.6667
$ws.Range("I103").Value = 600.4545000000001
$ws.Range("J103").Value = 1100
$ws.Range("K103").Value = 1801.3635
$ws.Range("L103").Value = 3300
$ws.Range("M103").Value = -922.3635000000002
$ws.Range("N103").Value = -5058

$ws.Range("H107").Value = 573.1111
$ws.Range("I107").Value = 278.33334
$ws.Range("J107").Value = 1162.6666
$ws.Range("K107").Value = 835.0000200000001
$ws.Range("L107").Value = 3487.9998
$ws.Range("M107").Value = 1084.99998
$ws.Range("N107").Value = -7327.9998

$ws.Range("H113").Value = 791.26
$ws.Range("I113").Value = 513.9
$ws.Range("J113").Value = 822.07776
$ws.Range("K113").Value = 1541.7
$ws.Range("L113").Value = 2466.23328
$ws.Range("M113").Value = 628.3000000000002
$ws.Range("N113").Value = -6806.23328

$ws.Range("H131").Value = 10753501
$ws.Range("I131").Value = 19231542
$ws.Range("J131").Value = 7463515
$ws.Range("K131").Value = 57694626
$ws.Range("L131").Value = 22390545
$ws.Range("M131").Value = -57689586
$ws.Range("N131").Value = -22400625


$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 11486006
$ws.Range("I70").Value = 17583912
$ws.Range("J70").Value = 7595.1763
$ws.Range("K70").Value = 17583912
$ws.Range("L70").Value = 7595.1763
$ws.Range("M70").Value = -17583642
$ws.Range("N70").Value = -8135.1763

$ws.Range("H73").Value = 11486006
$ws.Range("I73").Value = 17583912
$ws.Range("J73").Value = 7595.1763
$ws.Range("K73").Value = 17583912
$ws.Range("L73").Value = 7595.1763
$ws.Range("M73").Value = -17582976
$ws.Range("N73").Value = -9467.176299999999


$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H132").Value = 5344.4546
$ws.Range("I132").Value = 6229.4287
$ws.Range("J132").Value = 3795.75
$ws.Range("K132").Value = 18688.2861
$ws.Range("L132").Value = 11387.25
$ws.Range("M132").Value = -16158.2861
$ws.Range("N132").Value = -16447.25


$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H132").Value = 32612978
$ws.Range("I132").Value = 55557436
$ws.Range("J132").Value = 7694.6313
$ws.Range("K132").Value = 166672308
$ws.Range("L132").Value = 23083.8939
$ws.Range("M132").Value = -166669778
$ws.Range("N132").Value = -28143.8939

